# Rename "Cross references" to "Database references" (wc_lang.core.CrossReference ->
# wc_lang.core.DatabaseReference rename) and make it the active/selected sheet, matching
# the author's interactive session that ended with this sheet's tab selected.

$wb = $excel.ActiveWorkbook

$sheet = $wb.Worksheets.Item("Cross references")
$sheet.Name = "Database references"

# Activating the sheet updates workbookView.activeTab and swaps tabSelected="1"
# from the previously active sheet onto this one.
$sheet.Activate()
